$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data rows (16-22): full bond-length / angle records pulled from the
# HBIL-dft and S66x8 MP2/aVQZ datasets.
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "EtNH3-AcO-p1"
$ws.Range("C16").Value = 1.6373
$ws.Range("D16").Value = 163.11199999999999
$ws.Range("E16").Value = 1.6370199999999999
$ws.Range("F16").Value = 126.49
$ws.Range("G16").Value = 136.613

$ws.Range("B17").Value = "EtNH3-CF3SO3-p1"
$ws.Range("C17").Value = 1.8427800000000001
$ws.Range("D17").Value = 144.761
$ws.Range("E17").Value = 1.9504999999999999
$ws.Range("F17").Value = 109.34399999999999
$ws.Range("G17").Value = -16.667999999999999
$ws.Range("H17").Value = "135.733/110.716 bond length/angle for other H-bond"

$ws.Range("B18").Value = "EtNH3-Cl-p1"
$ws.Range("C18").Value = 1.9028700000000001
$ws.Range("D18").Value = 169.20699999999999
$ws.Range("F18").Value = "NA"
$ws.Range("G18").Value = "NA"
$ws.Range("H18").Value = "Cl counts?"

$ws.Range("B19").Value = "EtNH3-NO3-p1"
$ws.Range("C19").Value = 1.54939
$ws.Range("D19").Value = 168.76400000000001
$ws.Range("F19").Value = 105.295
$ws.Range("G19").Value = -33.479999999999997
$ws.Range("H19").Value = "Y"

$ws.Range("B20").Value = "EtNH3-TFA-p1"
$ws.Range("C20").Value = 1.45394
$ws.Range("D20").Value = 166.97
$ws.Range("F20").Value = 102.895
$ws.Range("G20").Value = 23.524999999999999
$ws.Range("H20").Value = "Y"

$ws.Range("B21").Value = "EtNH3-mOSO3-p1"
$ws.Range("C21").Value = 1.66107
$ws.Range("D21").Value = 161.18
$ws.Range("F21").Value = 103.449
$ws.Range("G21").Value = -29.39

$ws.Range("B22").Value = "EtNH3-mSO3-p1"
$ws.Range("C22").Value = 1.80535
$ws.Range("D22").Value = 142.755
$ws.Range("F22").Value = 109.631
$ws.Range("G22").Value = -21.157
$ws.Range("H22").Value = "1.80536/142.749 (NHO)"

# ---------------------------------------------------------------------------
# New label-only rows (23-49): TMEA / mim / mpyr series identifiers staged
# for the next data pull.
# ---------------------------------------------------------------------------
$ws.Range("B23").Value = "TMEA-AcO-p1"
$ws.Range("B24").Value = "TMEA-CF3SO3-p1"
$ws.Range("B25").Value = "TMEA-Cl-p1"
$ws.Range("B26").Value = "TMEA-NO3-p1"
$ws.Range("B27").Value = "TMEA-TFA-p1"
$ws.Range("B28").Value = "TMEA-mOSO3-p1"
$ws.Range("B29").Value = "TMEA-mSO3-p1"
$ws.Range("B30").Value = "mim-AcO-p1"
$ws.Range("B31").Value = "mim-CF3SO3-p1"
$ws.Range("B32").Value = "mim-CF3SO3-p2"
$ws.Range("B33").Value = "mim-Cl-p1"
$ws.Range("B34").Value = "mim-Cl-p2"
$ws.Range("B35").Value = "mim-NO3-p1"
$ws.Range("B36").Value = "mim-NO3-p2"
$ws.Range("B37").Value = "mim-TFA-p1"
$ws.Range("B38").Value = "mim-TFA-p2"
$ws.Range("B39").Value = "mim-TFA-p3"
$ws.Range("B40").Value = "mim-mOSO3-p1"
$ws.Range("B41").Value = "mim-mSO3-p1"
$ws.Range("B42").Value = "mim-mSO3-p2"
$ws.Range("B43").Value = "mpyr-AcO-p1"
$ws.Range("B44").Value = "mpyr-CF3SO3-p1"
$ws.Range("B45").Value = "mpyr-Cl-p1"
$ws.Range("B46").Value = "mpyr-NO3-p1"
$ws.Range("B47").Value = "mpyr-TFA-p1"
$ws.Range("B48").Value = "mpyr-mOSO3-p1"
$ws.Range("B49").Value = "mpyr-mSO3-p1"

# ---------------------------------------------------------------------------
# Number formats: column C (bond length) gets a new 5-decimal format,
# column G (dihedral) picks up the existing 3-decimal format already used
# by column D. Row 9 is a label-only row, so it is skipped on both ranges
# to avoid creating stray blank formatted cells.
# ---------------------------------------------------------------------------
$ws.Range("C2:C8").NumberFormat = "0.00000"
$ws.Range("C10:C22").NumberFormat = "0.00000"

$ws.Range("G2:G8").NumberFormat = "0.000"
$ws.Range("G10:G22").NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------------------
$ws.Range("D23").Select()
